$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diff")

# Insert a new row above the current row 8 (Julia/Matlab b=100*rand row),
# shifting the existing rows 8-19 down to 9-20, and inheriting formatting
# from the row above (row 7) just like Excel's native "Insert" does.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new "dt = 5 min" case.
$ws.Range("A8").Value = "Matlab transport, dt = 5 min, K=600, b=100"
$ws.Range("B8").Value = -0.0000000160910000000000021

# Match the post-edit selection recorded in the workbook.
$ws.Range("B5:B8").Select()
